# Automatic update of files.
#
# Applies the following data changes to the "Artfynd" sheet:
#   - B2, B3: Taxonsorteringsordning 79244 -> 79245
#   - Row 4 and Row 5 swap their A (Id), Q (Ost), R (Nord), Z (Starttid) and
#     AB (Sluttid) values
#   - B4, B5: Taxonsorteringsordning 79244 -> 79245
#   - B6: Taxonsorteringsordning 91829 -> 91830

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simple +1 increments on the Taxonsorteringsordning column.
$ws.Range("B2").Value = 79245
$ws.Range("B3").Value = 79245
$ws.Range("B4").Value = 79245
$ws.Range("B5").Value = 79245
$ws.Range("B6").Value = 91830

# Swap the row 4 / row 5 observation-specific values (Id, Ost, Nord,
# Starttid, Sluttid) - these two records traded places.
# Note: use .Value2 for reads ( .Value's getter is unreliable in this
# host), .Value for writes.
$row4Id = $ws.Range("A4").Value2
$row5Id = $ws.Range("A5").Value2
$ws.Range("A4").Value = $row5Id
$ws.Range("A5").Value = $row4Id

$row4Ost = $ws.Range("Q4").Value2
$row5Ost = $ws.Range("Q5").Value2
$ws.Range("Q4").Value = $row5Ost
$ws.Range("Q5").Value = $row4Ost

$row4Nord = $ws.Range("R4").Value2
$row5Nord = $ws.Range("R5").Value2
$ws.Range("R4").Value = $row5Nord
$ws.Range("R5").Value = $row4Nord

$row4Start = $ws.Range("Z4").Value2
$row5Start = $ws.Range("Z5").Value2
$ws.Range("Z4").Value = $row5Start
$ws.Range("Z5").Value = $row4Start

$row4End = $ws.Range("AB4").Value2
$row5End = $ws.Range("AB5").Value2
$ws.Range("AB4").Value = $row5End
$ws.Range("AB5").Value = $row4End
